$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.747.75"
$ws.Range("E2").Value = "  -1.39%  "

$ws.Range("D3").Value = "3.409.07"
$ws.Range("E3").Value = "  -0.26%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.87"
$ws.Range("E5").Value = "  -0.70%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.69"
$ws.Range("E6").Value = "  +0.32%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "3.411.48"
$ws.Range("E8").Value = "  -0.30%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.551"
$ws.Range("E9").Value = "  -9.30%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.28"
$ws.Range("E10").Value = "  +2.05%  "

$ws.Range("E11").Value = "  -3.10%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.422"
$ws.Range("E12").Value = "  -4.01%  "

$ws.Range("D13").Value = "4.001.85"
$ws.Range("E13").Value = "  -0.14%  "

$ws.Range("E14").Value = "  +0.77%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.02"
$ws.Range("E15").Value = "  -1.98%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000172"
$ws.Range("E16").Value = "  -8.47%  "

$ws.Range("D17").Value = "63.866.43"
$ws.Range("E17").Value = "  -1.21%  "

$ws.Range("D18").Value = "3.377.65"
$ws.Range("E18").Value = "  -1.12%  "

$ws.Range("E19").Value = "  -4.37%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.52"
$ws.Range("E20").Value = "  -2.34%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "376.12"
$ws.Range("E21").Value = "  -1.03%  "

$ws.Range("E22").Value = "  -1.99%  "

$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.47"
$ws.Range("E24").Value = "  -0.99%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.514"
$ws.Range("E25").Value = "  -6.29%  "

$ws.Range("E26").Value = "  -2.69%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.55"
$ws.Range("E27").Value = "  -5.67%  "

$ws.Range("E28").Value = "  -0.66%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.11%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.97"
$ws.Range("E30").Value = "  -3.23%  "

$ws.Range("E31").Value = "  -6.09%  "

$ws.Range("E32").Value = "  -0.69%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "22.81"
$ws.Range("E33").Value = "  -1.60%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.02"
$ws.Range("E34").Value = "  -0.77%  "

$ws.Range("E35").Value = "  -5.83%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "159.43"
$ws.Range("E36").Value = "  -0.60%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.831"
$ws.Range("E37").Value = "  +7.84%  "

$ws.Range("E38").Value = "  -5.84%  "

$ws.Range("D39").Value = "2.813.17"
$ws.Range("E39").Value = "  -2.45%  "

$ws.Range("E40").Value = "  -3.64%  "

$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "42.89"
$ws.Range("E41").Value = "  -0.22%  "

$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "25.70"
$ws.Range("E42").Value = "  -2.88%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.46"
$ws.Range("E43").Value = "  -5.07%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.42"
$ws.Range("E44").Value = "  -3.64%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.85"
$ws.Range("E45").Value = "  +0.50%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0302"
$ws.Range("E46").Value = "  -4.10%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "335.58"
$ws.Range("E47").Value = "  +4.38%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.36"
$ws.Range("E48").Value = "  +6.53%  "

$ws.Range("E49").Value = "  -2.32%  "

$ws.Range("E50").Value = "  -4.66%  "

$ws.Range("E51").Value = "  -3.91%  "
